# Apply the "Add files via upload" edit to nameList.xlsx (Sheet1)
#
# Summary of the change:
#  - Column A gets two new names inserted into the list ("Yvoone" before
#    "Lam Kuen", and "JonathanW" right after "Jonathan"), pushing the
#    remaining names down.
#  - Column B loses "La Son" (names shift up) and gains two new names
#    at the bottom ("Somingtat", "SomingtatW").
#  - Column C loses "Billy" from its original spot (it is relocated to a
#    new row 15) and gains a new name at the bottom ("Fion").
#  - The used range grows from A1:C12 to A1:C15.
#  - The active selection moves to B12.
#  - The sheet gets an explicit A4 / portrait page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Column A -------------------------------------------------------
$ws.Range("A2").Value2  = "Yvoone"
$ws.Range("A3").Value2  = "Lam Kuen"
$ws.Range("A4").Value2  = "Leo"
$ws.Range("A5").Value2  = "Francis"
$ws.Range("A6").Value2  = "Small Tony"
$ws.Range("A7").Value2  = "Anthony "
$ws.Range("A8").Value2  = "Jonathan"
$ws.Range("A9").Value2  = "JonathanW"
$ws.Range("A10").Value2 = "Rohda"
$ws.Range("A11").Value2 = "Patrick"
$ws.Range("A12").Value2 = "Faye"

# ---- Column B -------------------------------------------------------
$ws.Range("B7").Value2  = "Begger"
$ws.Range("B8").Value2  = "Harry Cheung"
$ws.Range("B9").Value2  = "See Fu"
$ws.Range("B10").Value2 = "Somingtat"
$ws.Range("B11").Value2 = "SomingtatW"

# ---- Column C -------------------------------------------------------
$ws.Range("C4").Value2  = "Tim"
$ws.Range("C5").Value2  = "Arun"
$ws.Range("C6").Value2  = "Natalie"
$ws.Range("C7").Value2  = "Denis"
$ws.Range("C8").Value2  = "Gabriel"
$ws.Range("C9").Value2  = "Hao"
$ws.Range("C10").Value2 = "Arbinnav"
$ws.Range("C11").Value2 = "Denvendra"
$ws.Range("C12").Value2 = "Fion"
$ws.Range("C15").Value2 = "Billy"

# ---- Page setup (A4, portrait) --------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# ---- Selection --------------------------------------------------------
[void]$ws.Range("B12").Select()
